$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values - automatic update of electricity spot prices (day shifted +1)
$ws.Range("A2").Value = 45969
$ws.Range("B2").Value = 67.59
$ws.Range("C2").Value = 56.35
$ws.Range("D2").Value = 50.98
$ws.Range("E2").Value = 45.53
$ws.Range("F2").Value = 35.71
$ws.Range("G2").Value = 36.99
$ws.Range("H2").Value = 47.35
$ws.Range("I2").Value = 50.9
$ws.Range("J2").Value = 38.51
$ws.Range("K2").Value = 4.42
$ws.Range("L2").Value = 3.53
$ws.Range("M2").Value = 9.11
$ws.Range("N2").Value = 16.59
$ws.Range("O2").Value = 22.4
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 65.27
$ws.Range("T2").Value = 87.33
$ws.Range("U2").Value = 88.84
$ws.Range("V2").Value = 84.11
$ws.Range("W2").Value = 78.27
$ws.Range("X2").Value = 75.82
$ws.Range("Y2").Value = 73.01
$ws.Range("Z2").Value = 44.73
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 77.8
$ws.Range("AD2").Value = 88.08
$ws.Range("AF2").Value = 81.19
$ws.Range("AG2").Value = "4h-16h"
